# The document ends with three paragraphs that need to collapse into one:
#   ...,  "，，，，，，",  (empty),  "。。。。" [+ the _GoBack bookmark]
# Target: a single empty paragraph that keeps the first paragraph's
# formatting and the bookmark from the last paragraph.
$d = $word.ActiveDocument

# 1) Strip the placeholder text out of the two runs, leaving three
#    paragraphs that are now empty (but still separate).
$d.Content.Find.Execute("，，，，，，", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
$d.Content.Find.Execute("。。。。", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 2) Collapse the three trailing empty paragraphs into one by deleting the
#    paragraph marks between them. The last paragraph (with the _GoBack
#    bookmark) is always $d.Paragraphs.Count; the first of the three is
#    three paragraphs before that. Deleting that paragraph's Range (which
#    spans through its own end-of-paragraph mark) merges it with the next
#    paragraph; repeating once more merges the result with the final one.
$firstIndex = $d.Paragraphs.Count - 2

$p = $d.Paragraphs.Item($firstIndex)
$d.Range($p.Range.Start, $p.Range.End).Delete() | Out-Null

$p = $d.Paragraphs.Item($firstIndex)
$d.Range($p.Range.Start, $p.Range.End).Delete() | Out-Null
